$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values (A1, B1, C1 on the single sheet)
$ws.Range("A1").Value = 149.06146062471865
$ws.Range("B1").Value = 4.7356777681544777
$ws.Range("C1").Value = 0.81663405088062624

# Update column widths for columns B and C.
# Excel's ColumnWidth property is expressed in "characters" and gets
# quantized to the screen pixel grid when stored back to the OOXML
# <col width="..."/> attribute (raw = ColumnWidth + 5/MaxDigitWidth).
# Choose ColumnWidth values whose resulting raw width is the closest
# achievable value to the target widths (11.7109375 and 12.7109375).
$ws.Columns.Item(2).ColumnWidth = 10.83
$ws.Columns.Item(3).ColumnWidth = 11.83
